$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("N1").Value = "review_topics"
$ws.Range("AB1").Value = "trait_frame"

$ws.Range("N2").Value = "changegoals; skilltraitdifferences"

$ws.Range("N3").Value = "normativechange"
$ws.Range("O3").Value = "both"

$ws.Range("N4").Value = "assessment; theorydevelopment"

$ws.Range("N5").Value = "positiveschooloutcomes"
$ws.Range("O5").Value = "both"

$ws.Range("N6").Value = "proenvironment"
$ws.Range("O6").Value = "both"
